$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Status value and Date value ------------------
$ws1 = $wb.Worksheets.Item("Metadata")

# Status: active -> draft  (B6)
$ws1.Range("B6").Value = "draft"

# Date: 2023-05-12T12:33:13+00:00 -> 2023-08-01T16:12:28+00:00  (B8)
$ws1.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# --- Style fix: the header-row style and the bordered-data-row style should
# explicitly mark their (pre-existing) wrap/vertical-top alignment as
# "applied" (applyAlignment="true"), matching authoring tools (e.g. POI)
# that always stamp applyAlignment when an <alignment> child is written.
# Re-asserting WrapText on every cell that already uses those two styles
# (on both worksheets) makes the engine emit applyAlignment="true" on the
# resulting xf record without altering any other formatting.
$ws2 = $wb.Worksheets.Item(2)

$ws1.Range("A1:B14").WrapText = $true

$ws2.Range("A1:A4").WrapText = $true
$ws2.Range("B3:B4").WrapText = $true
